$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = "NK Sesvete"
$ws.Cells.Item(5, 7).Value = "HNK Cibalia"
$ws.Cells.Item(9, 2).Value = 6834733
$ws.Cells.Item(9, 7).Value = "NK Croatia Zmijavci"
$ws.Cells.Item(9, 8).Value = 1
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = "H"
$ws.Cells.Item(9, 11).Value = 1.65
$ws.Cells.Item(9, 12).Value = 3.5
$ws.Cells.Item(9, 13).Value = 4.5
$ws.Cells.Item(9, 14).Value = 1.909
$ws.Cells.Item(9, 15).Value = 3.3
$ws.Cells.Item(9, 16).Value = 3.3
$ws.Cells.Item(9, 17).Value = -0.5
$ws.Cells.Item(9, 18).Value = 2
$ws.Cells.Item(9, 19).Value = 1.8
$ws.Cells.Item(9, 20).Value = 2.25
$ws.Cells.Item(9, 21).Value = 1.95
$ws.Cells.Item(9, 22).Value = 1.85
$ws.Cells.Item(9, 23).Value = 0.909
$ws.Cells.Item(9, 25).Value = -1
$ws.Cells.Item(9, 26).Value = 1
$ws.Cells.Item(9, 27).Value = -1
$ws.Cells.Item(9, 29).Value = 0.8500000000000001
$ws.Cells.Item(10, 2).Value = 6834729
$ws.Cells.Item(10, 7).Value = "Bijelo Brdo"
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 10).Value = "D"
$ws.Cells.Item(10, 11).Value = 2.1
$ws.Cells.Item(10, 12).Value = 3.2
$ws.Cells.Item(10, 13).Value = 3.1
$ws.Cells.Item(10, 14).Value = 2.05
$ws.Cells.Item(10, 15).Value = 3.25
$ws.Cells.Item(10, 16).Value = 3.25
$ws.Cells.Item(10, 17).Value = -0.25
$ws.Cells.Item(10, 18).Value = 1.8
$ws.Cells.Item(10, 19).Value = 2
$ws.Cells.Item(10, 20).Value = 2.5
$ws.Cells.Item(10, 21).Value = 2
$ws.Cells.Item(10, 22).Value = 1.8
$ws.Cells.Item(10, 23).Value = -1
$ws.Cells.Item(10, 24).Value = 2.25
$ws.Cells.Item(10, 26).Value = -0.5
$ws.Cells.Item(10, 27).Value = 0.5
$ws.Cells.Item(10, 29).Value = 0.8
$ws.Cells.Item(11, 2).Value = 7045999
$ws.Cells.Item(11, 7).Value = "NK Dugopolje"
$ws.Cells.Item(11, 9).Value = 1
$ws.Cells.Item(11, 10).Value = "A"
$ws.Cells.Item(11, 11).Value = 2.4
$ws.Cells.Item(11, 12).Value = 3.25
$ws.Cells.Item(11, 13).Value = 2.6
$ws.Cells.Item(11, 14).Value = 2.4
$ws.Cells.Item(11, 16).Value = 2.6
$ws.Cells.Item(11, 17).Value = 0
$ws.Cells.Item(11, 18).Value = 1.775
$ws.Cells.Item(11, 19).Value = 2.025
$ws.Cells.Item(11, 20).Value = 2
$ws.Cells.Item(11, 21).Value = 1.775
$ws.Cells.Item(11, 22).Value = 2.025
$ws.Cells.Item(11, 24).Value = -1
$ws.Cells.Item(11, 25).Value = 1.6
$ws.Cells.Item(11, 26).Value = -1
$ws.Cells.Item(11, 27).Value = 1.025
$ws.Cells.Item(11, 29).Value = 1.025
$ws.Cells.Item(15, 7).Value = "NK Sesvete"
$ws.Cells.Item(17, 6).Value = "HNK Cibalia"
$ws.Cells.Item(18, 7).Value = "NK Solin"
$ws.Cells.Item(22, 6).Value = "NK Sesvete"
$ws.Cells.Item(23, 7).Value = "HNK Cibalia"
$ws.Cells.Item(25, 6).Value = "NK Solin"
$ws.Cells.Item(27, 7).Value = "NK Sesvete"
$ws.Cells.Item(28, 6).Value = "HNK Cibalia"
$ws.Cells.Item(29, 7).Value = "NK Solin"
$ws.Cells.Item(32, 6).Value = "NK Solin"
$ws.Cells.Item(36, 6).Value = "NK Sesvete"
$ws.Cells.Item(36, 7).Value = "HNK Cibalia"
$ws.Cells.Item(41, 6).Value = "HNK Cibalia"
$ws.Cells.Item(41, 7).Value = "NK Solin"
$ws.Cells.Item(43, 7).Value = "NK Sesvete"
$ws.Cells.Item(44, 7).Value = "HNK Cibalia"
$ws.Cells.Item(45, 6).Value = "NK Solin"
$ws.Cells.Item(47, 6).Value = "NK Sesvete"
$ws.Cells.Item(50, 7).Value = "NK Solin"
$ws.Cells.Item(51, 7).Value = "NK Sesvete"
$ws.Cells.Item(53, 6).Value = "HNK Cibalia"
$ws.Cells.Item(56, 6).Value = "NK Solin"
$ws.Cells.Item(59, 6).Value = "NK Sesvete"
$ws.Cells.Item(61, 7).Value = "HNK Cibalia"
$ws.Cells.Item(63, 6).Value = "HNK Cibalia"
$ws.Cells.Item(64, 6).Value = "NK Sesvete"
$ws.Cells.Item(64, 7).Value = "NK Solin"
$ws.Cells.Item(68, 6).Value = "NK Solin"
$ws.Cells.Item(70, 7).Value = "HNK Cibalia"
$ws.Cells.Item(71, 7).Value = "NK Sesvete"
$ws.Cells.Item(76, 7).Value = "HNK Cibalia"
$ws.Cells.Item(77, 7).Value = "NK Solin"
$ws.Cells.Item(78, 7).Value = "NK Sesvete"
$ws.Cells.Item(80, 6).Value = "NK Solin"
$ws.Cells.Item(84, 6).Value = "NK Sesvete"
$ws.Cells.Item(85, 7).Value = "HNK Cibalia"
$ws.Cells.Item(86, 6).Value = "HNK Cibalia"
$ws.Cells.Item(87, 7).Value = "NK Solin"
$ws.Cells.Item(91, 7).Value = "NK Sesvete"
$ws.Cells.Item(92, 6).Value = "NK Sesvete"
$ws.Cells.Item(93, 6).Value = "NK Solin"
$ws.Cells.Item(95, 7).Value = "HNK Cibalia"
$ws.Cells.Item(98, 7).Value = "NK Solin"
$ws.Cells.Item(99, 6).Value = "HNK Cibalia"
$ws.Cells.Item(99, 7).Value = "NK Sesvete"
$ws.Cells.Item(105, 6).Value = "NK Solin"
$ws.Cells.Item(105, 7).Value = "HNK Cibalia"
$ws.Cells.Item(106, 6).Value = "NK Sesvete"
$ws.Cells.Item(111, 7).Value = "NK Sesvete"
$ws.Cells.Item(112, 6).Value = "HNK Cibalia"
$ws.Cells.Item(114, 7).Value = "NK Solin"
$ws.Cells.Item(116, 6).Value = "NK Solin"
$ws.Cells.Item(119, 6).Value = "NK Sesvete"
$ws.Cells.Item(121, 7).Value = "HNK Cibalia"
$ws.Cells.Item(124, 7).Value = "NK Sesvete"
$ws.Cells.Item(125, 7).Value = "NK Solin"
$ws.Cells.Item(127, 6).Value = "HNK Cibalia"
$ws.Cells.Item(128, 6).Value = "NK Solin"
$ws.Cells.Item(128, 7).Value = "NK Sesvete"
$ws.Cells.Item(130, 2).Value = 6834845
$ws.Cells.Item(130, 6).Value = "NK Dubrava Zagreb"
$ws.Cells.Item(130, 7).Value = "HNK Sibenik"
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 9).Value = 3
$ws.Cells.Item(130, 10).Value = "A"
$ws.Cells.Item(130, 11).Value = 3.75
$ws.Cells.Item(130, 12).Value = 3.5
$ws.Cells.Item(130, 13).Value = 1.8
$ws.Cells.Item(130, 14).Value = 4.5
$ws.Cells.Item(130, 15).Value = 3.6
$ws.Cells.Item(130, 16).Value = 1.65
$ws.Cells.Item(130, 17).Value = 0.75
$ws.Cells.Item(130, 18).Value = 1.925
$ws.Cells.Item(130, 19).Value = 1.875
$ws.Cells.Item(130, 20).Value = 2.25
$ws.Cells.Item(130, 21).Value = 1.9
$ws.Cells.Item(130, 22).Value = 1.9
$ws.Cells.Item(130, 23).Value = -1
$ws.Cells.Item(130, 25).Value = 0.6499999999999999
$ws.Cells.Item(130, 26).Value = -1
$ws.Cells.Item(130, 27).Value = 0.875
$ws.Cells.Item(130, 28).Value = 0.8999999999999999
$ws.Cells.Item(130, 29).Value = -1
$ws.Cells.Item(131, 2).Value = 6834844
$ws.Cells.Item(131, 6).Value = "NK Jarun"
$ws.Cells.Item(131, 7).Value = "HNK Cibalia"
$ws.Cells.Item(131, 8).Value = 2
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = "H"
$ws.Cells.Item(131, 11).Value = 1.727
$ws.Cells.Item(131, 12).Value = 3.6
$ws.Cells.Item(131, 13).Value = 4
$ws.Cells.Item(131, 14).Value = 2.05
$ws.Cells.Item(131, 15).Value = 3.4
$ws.Cells.Item(131, 16).Value = 3.1
$ws.Cells.Item(131, 17).Value = -0.25
$ws.Cells.Item(131, 18).Value = 1.85
$ws.Cells.Item(131, 19).Value = 1.95
$ws.Cells.Item(131, 20).Value = 2.5
$ws.Cells.Item(131, 21).Value = 1.925
$ws.Cells.Item(131, 22).Value = 1.875
$ws.Cells.Item(131, 23).Value = 1.05
$ws.Cells.Item(131, 25).Value = -1
$ws.Cells.Item(131, 26).Value = 0.8500000000000001
$ws.Cells.Item(131, 27).Value = -1
$ws.Cells.Item(131, 28).Value = -1
$ws.Cells.Item(131, 29).Value = 0.875
$ws.Cells.Item(134, 6).Value = "NK Sesvete"
$ws.Cells.Item(136, 7).Value = "NK Solin"
$ws.Cells.Item(139, 6).Value = "HNK Cibalia"
$ws.Cells.Item(140, 6).Value = "NK Solin"
$ws.Cells.Item(144, 7).Value = "HNK Cibalia"
$ws.Cells.Item(145, 7).Value = "NK Sesvete"
$ws.Cells.Item(148, 7).Value = "NK Solin"
$ws.Cells.Item(149, 6).Value = "NK Sesvete"
$ws.Cells.Item(150, 6).Value = "HNK Cibalia"
$ws.Cells.Item(153, 6).Value = "NK Solin"
$ws.Cells.Item(153, 7).Value = "HNK Cibalia"
$ws.Cells.Item(155, 7).Value = "NK Sesvete"
